$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 110947080
$ws.Range("B2").Value = 96348
$ws.Range("Q2").Value = 614510.1279625499
$ws.Range("R2").Value = 6657642.00361704
$ws.Range("Z2").Value = '13:11'
$ws.Range("AB2").Value = '13:11'

# Row 3
$ws.Range("A3").Value = 110947619
$ws.Range("B3").Value = 96348
$ws.Range("Q3").Value = 614481.0570550568
$ws.Range("R3").Value = 6657755.583492418
$ws.Range("Z3").Value = '13:36'
$ws.Range("AB3").Value = '13:36'

# Row 4
$ws.Range("A4").Value = 110947019
$ws.Range("B4").Value = 103288
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 221144
$ws.Range("F4").Value = 'Grönpyrola'
$ws.Range("G4").Value = 'Pyrola chlorantha'
$ws.Range("H4").Value = 'Sw.'
$ws.Range("Q4").Value = 614534.747918217
$ws.Range("R4").Value = 6657623.271711768
$ws.Range("Z4").Value = '13:07'
$ws.Range("AB4").Value = '13:07'

# Row 5
$ws.Range("A5").Value = 110947410
$ws.Range("B5").Value = 96348
$ws.Range("Q5").Value = 614542.6751059515
$ws.Range("R5").Value = 6657706.507382731
$ws.Range("Z5").Value = '13:28'
$ws.Range("AB5").Value = '13:28'

# Row 6
$ws.Range("A6").Value = 110947351
$ws.Range("B6").Value = 89369
$ws.Range("E6").Value = 5447
$ws.Range("F6").Value = 'Vedticka'
$ws.Range("G6").Value = 'Fuscoporia viticola'
$ws.Range("H6").Value = '(Schwein.) Murrill'
$ws.Range("Q6").Value = 614544.7041997212
$ws.Range("R6").Value = 6657689.572886499
$ws.Range("Z6").Value = '13:23'
$ws.Range("AB6").Value = '13:23'

# Row 7
$ws.Range("A7").Value = 110948255
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = 'Knärot'
$ws.Range("G7").Value = 'Goodyera repens'
$ws.Range("H7").Value = '(L.) R. Br.'
$ws.Range("Q7").Value = 614423.4236839975
$ws.Range("R7").Value = 6657789.286310961
$ws.Range("Z7").Value = '14:02'
$ws.Range("AB7").Value = '14:02'

# Row 8
$ws.Range("A8").Value = 110948236
$ws.Range("B8").Value = 78604
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 6461
$ws.Range("F8").Value = 'Norrlandslav'
$ws.Range("G8").Value = 'Nephroma arcticum'
$ws.Range("H8").Value = '(L.) Torss.'
$ws.Range("Q8").Value = 614411.6690967374
$ws.Range("R8").Value = 6657796.919702402
$ws.Range("Z8").Value = '13:58'
$ws.Range("AB8").Value = '13:58'

# Row 9
$ws.Range("A9").Value = 110947035
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = 'Knärot'
$ws.Range("G9").Value = 'Goodyera repens'
$ws.Range("H9").Value = '(L.) R. Br.'
$ws.Range("Q9").Value = 614524.5089896603
$ws.Range("R9").Value = 6657630.452345544
$ws.Range("Z9").Value = '13:05'
$ws.Range("AB9").Value = '13:05'

# Row 10
$ws.Range("A10").Value = 110947491
$ws.Range("B10").Value = 96348
$ws.Range("Q10").Value = 614553.3115441641
$ws.Range("R10").Value = 6657734.834170708
$ws.Range("Z10").Value = '13:31'
$ws.Range("AB10").Value = '13:31'

# Row 11
$ws.Range("A11").Value = 110948416
$ws.Range("B11").Value = 96348
$ws.Range("P11").Value = 'Björkmossen (Björkmossen), Upl'
$ws.Range("Q11").Value = 614392.0112977019
$ws.Range("R11").Value = 6657768.813506908
$ws.Range("Z11").Value = '14:11'
$ws.Range("AB11").Value = '14:11'

# Row 12
$ws.Range("A12").Value = 110948582
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = 'VU'
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = 'Knärot'
$ws.Range("G12").Value = 'Goodyera repens'
$ws.Range("H12").Value = '(L.) R. Br.'
$ws.Range("P12").Value = 'Bredmossen (Bredmossen), Upl'
$ws.Range("Q12").Value = 614461.1684684249
$ws.Range("R12").Value = 6657815.455187102
$ws.Range("Z12").Value = '14:23'
$ws.Range("AB12").Value = '14:23'

# Row 13
$ws.Range("A13").Value = 110947380
$ws.Range("B13").Value = 96265
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 219790
$ws.Range("F13").Value = 'Fläcknycklar'
$ws.Range("G13").Value = 'Dactylorhiza maculata'
$ws.Range("H13").Value = '(L.) Soó'
$ws.Range("Q13").Value = 614569.4543792737
$ws.Range("R13").Value = 6657698.841700966
$ws.Range("Z13").Value = '13:25'
$ws.Range("AB13").Value = '13:25'

# Row 14
$ws.Range("B14").Value = 96348
